# Fix Training Data Issue: the BF column ("Date") held values in the wrong
# format (e.g. "4-24-2013-14") because of how the NBA stats site showed
# dates spanning two season years. Correct it to ISO "yyyy-mm-dd" format
# (e.g. "2014-04-24") for every data row, keeping the values as plain text
# (not auto-converted to an Excel date serial).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF" + $row)
    $old = [string]$cell.Value2
    if ($old -eq "4-24-2013-14") {
        # Leading apostrophe forces text entry so Excel doesn't reinterpret
        # the ISO-formatted string as a date value.
        $cell.Value = "'2014-04-24"
        # Re-apply the default "Normal" style so the cell keeps its original
        # (unstyled) formatting instead of the quote-prefixed text format.
        $cell.Style = "Normal"
    }
}
